$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that currently sits in the
#    empty paragraph above "Important note about app". It gets
#    recreated later, inside the dropbox hyperlink run.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) "Important note about app" -> bold + underlined, and append
#    " update" as its own run (still bold/underlined) so the
#    paragraph reads "Important note about app update".
# ------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("Important note about app", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$noteStart = $find.Start
$noteEnd = $find.End

$noteXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>Important note about app</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> update</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$noteRange = $d.Range($noteStart, $noteEnd)
$noteRange.InsertXML($noteXml)

# ------------------------------------------------------------------
# 3) Split the dropbox hyperlink run into three runs
#    ("https://w" / "w" / "w.dropbox.com/...") and drop the
#    _GoBack bookmark right after the first piece, matching the
#    cursor-split that Word leaves behind.
# ------------------------------------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("https://www.dropbox.com/sh/vtx6425dnxce2fl/ijrUd_wzHb", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$linkStart = $find2.Start
$linkEnd = $find2.End

$linkXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:hyperlink r:id="rId6" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>https://w</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>w</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>w.dropbox.com/sh/vtx6425dnxce2fl/ijrUd_wzHb</w:t></w:r></w:hyperlink></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$linkRange = $d.Range($linkStart, $linkEnd)
$linkRange.InsertXML($linkXml)

# ------------------------------------------------------------------
# 4) Register the (built-in, previously unused) FollowedHyperlink
#    character style so it shows up explicitly in styles.xml /
#    stylesWithEffects.xml.
# ------------------------------------------------------------------
if (-not $d.Styles.Exists("FollowedHyperlink")) {
    $d.Styles.Add("FollowedHyperlink", 2)
}
